# Table_S1_Sylvatic_DENV-2_Cynomolgus_Macaques.xlsx
# Commit message: "Changed Table names (invert S3 and S4)"
#
# Semantic changes applied:
#   1. Rename the first worksheet from "Trade-Off_Project_Cyno_Inf_Sylv" to
#      "Data".
#   2. Make that worksheet ("Data") the active / selected sheet - it used to
#      be "Column Information" that was active/selected.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item(1)
$wsInfo = $wb.Worksheets.Item(2)

# --- 1. Rename the first sheet ------------------------------------------
$wsData.Name = "Data"

# --- 2. Make it the active sheet / tab -----------------------------------
# Previously "Column Information" (sheet 2) was the active tab; the edit
# flips this so the renamed "Data" sheet (sheet 1) becomes active.
$wsData.Activate()

# --- 3. Tidy up a handful of redundant (duplicate) cell formats on the
#        "Column Information" sheet. These cells were pointing at a cell
#        style that is a byte-for-byte duplicate of another, already used,
#        cell style; re-applying the equivalent "General" number format
#        lets the engine collapse them onto the de-duplicated style while
#        leaving every value and every other piece of formatting (bold
#        rows, Calibri rows, column widths, etc.) untouched.
$wsInfo.Range("A1:B1").NumberFormat = "General"
$wsInfo.Range("A4:B9").NumberFormat = "General"
$wsInfo.Range("A10:A14").NumberFormat = "General"
$wsInfo.Range("A15:B52").NumberFormat = "General"
$wsInfo.Range("A53:A62").NumberFormat = "General"
$wsInfo.Range("A63:B109").NumberFormat = "General"
